$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D6", "D7", "D11", "D12", "D14", "D16", "D18", "D19", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D37", "D38", "D39", "D41", "D43", "D45", "D46", "D47", "D50", "D51")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "64.222.10"
$ws.Range("E2").Value = "  -3.26%  "
$ws.Range("D3").Value = "3.345.19"
$ws.Range("E3").Value = "  -5.13%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "553.55"
$ws.Range("E5").Value = "  -5.19%  "
$ws.Range("D6").Value = "174.22"
$ws.Range("E6").Value = "  -3.10%  "
$ws.Range("D7").Value = "0.614"
$ws.Range("E7").Value = "  -2.62%  "
$ws.Range("D8").Value = "3.337.95"
$ws.Range("E8").Value = "  -5.17%  "
$ws.Range("E10").Value = "  -2.25%  "
$ws.Range("D11").Value = "0.161"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").Value = "53.75"
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("D14").Value = "9.03"
$ws.Range("E14").Value = "  -3.05%  "
$ws.Range("D15").Value = "3.875.17"
$ws.Range("D16").Value = "18.31"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.353.54"
$ws.Range("E17").Value = "  -5.02%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.117"
$ws.Range("E18").Value = "  -3.24%  "
$ws.Range("D19").Value = "11.77"
$ws.Range("E19").Value = "  -2.74%  "
$ws.Range("D20").Value = "64.107.73"
$ws.Range("E20").Value = "  -3.38%  "
$ws.Range("D21").Value = "0.976"
$ws.Range("E21").Value = "  -3.19%  "
$ws.Range("D22").Value = "424.04"
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("D23").Value = "4.84"
$ws.Range("E23").Value = "  +9.88%  "
$ws.Range("D24").Value = "4.09"
$ws.Range("E24").Value = "  -4.80%  "
$ws.Range("E25").Value = "  -2.08%  "
$ws.Range("D26").Value = "13.15"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "10.74"
$ws.Range("E27").Value = "  -4.12%  "
$ws.Range("D28").Value = "2.82"
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("D29").Value = "8.62"
$ws.Range("E29").Value = "  -5.77%  "
$ws.Range("D30").Value = "29.68"
$ws.Range("E30").Value = "  -2.64%  "
$ws.Range("D31").Value = "6.66"
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "11.42"
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "568.14"
$ws.Range("E33").Value = "  -5.96%  "
$ws.Range("E34").Value = "  -3.71%  "
$ws.Range("D35").Value = "58.14"
$ws.Range("E35").Value = "  -3.14%  "
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").Value = "0.142"
$ws.Range("E37").Value = "  -8.08%  "
$ws.Range("D38").Value = "3.48"
$ws.Range("E38").Value = "  -4.59%  "
$ws.Range("D39").Value = "35.54"
$ws.Range("E39").Value = "  -4.79%  "
$ws.Range("D40").Value = "0.0₃0753"
$ws.Range("E40").Value = "  -6.97%  "
$ws.Range("D41").Value = "0.366"
$ws.Range("E41").Value = "  -5.15%  "
$ws.Range("D42").Value = "3.096.09"
$ws.Range("E42").Value = "  -4.95%  "
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  -5.55%  "
$ws.Range("D45").Value = "3.22"
$ws.Range("E45").Value = "  -3.75%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "2.46"
$ws.Range("E46").Value = "  -3.92%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0406"
$ws.Range("E47").Value = "  -4.00%  "
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("E49").Value = "  -5.11%  "
$ws.Range("D50").Value = "136.44"
$ws.Range("E50").Value = "  -2.57%  "
$ws.Range("D51").Value = "8.21"
$ws.Range("E51").Value = "  -5.57%  "
